$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (pushes existing rows 4+ down to 5+)
$ws.Rows.Item(4).Insert()

# Update row 2 (Fragment) and row 3 (Column) "Notes" column text to the revised wording
$ws.Range("F2").Value = "Nest within <div type=""edition"">; must be followed by <ab>. Not necessary except when fragments are discrete and do not join?"
$ws.Range("F3").Value = "Nest within <div type=""edition"">; must be followed by <ab>"

# Populate the newly inserted row 4 describing the "Face (a)" text division
$ws.Range("B4").Value = "Text on separate surfaces of stone"
$ws.Range("C4").Value = "Face (a)"
$ws.Range("E4").Value = "<div type=""textpart"" subtype=""face"" n=""a"">"
$ws.Range("F4").Value = "Nest within <div type=""edition"">; must be followed by <ab>"

# Update the selected cell to match the saved view state
$ws.Range("F4").Select() | Out-Null
